# Apply updated odds/values to the "Jogos da Semana" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 1.8
$ws.Range("U2").Value = 1.36
$ws.Range("V2").Value = 3
$ws.Range("AA2").Value = 13

# Row 6
$ws.Range("H6").Value = 11
$ws.Range("I6").Value = 29
$ws.Range("W6").Value = 2.25
$ws.Range("X6").Value = 1.57
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 6.5
$ws.Range("AG6").Value = 41
$ws.Range("AH6").Value = 101
$ws.Range("AJ6").Value = 151
$ws.Range("AM6").Value = 201
$ws.Range("AN6").Value = 126

# Row 7
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 15
$ws.Range("Q7").Value = 1.67
$ws.Range("R7").Value = 2.15

# Row 8
$ws.Range("G8").Value = 1.4
$ws.Range("I8").Value = 6.5
$ws.Range("J8").Value = 1.91
$ws.Range("L8").Value = 6.5
$ws.Range("U8").Value = 1.3
$ws.Range("V8").Value = 3.4
$ws.Range("AJ8").Value = 34
